$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "28.934.83"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "1.881.51"
$ws.Range("E3").Value = "  -0.21%  "
Set-TextValue $ws.Range("D4") "1.002"
$ws.Range("E4").Value = "  -0.42%  "
Set-TextValue $ws.Range("D5") "325.13"
Set-TextValue $ws.Range("D6") "1.002"
$ws.Range("E6").Value = "  -0.43%  "
Set-TextValue $ws.Range("D7") "0.4612"
$ws.Range("E7").Value = "  +0.17%  "
Set-TextValue $ws.Range("D8") "0.3878"
$ws.Range("E8").Value = "  +0.17%  "
Set-TextValue $ws.Range("D9") "0.07849"
$ws.Range("E9").Value = "  -0.39%  "
Set-TextValue $ws.Range("D10") "0.9860"
$ws.Range("E10").Value = "  -1.69%  "
Set-TextValue $ws.Range("D11") "21.79"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").Value = "1.874.94"
$ws.Range("E12").Value = "  -2.64%  "
Set-TextValue $ws.Range("D13") "7.001"
$ws.Range("E13").Value = "  -1.12%  "
Set-TextValue $ws.Range("D14") "5.651"
$ws.Range("E14").Value = "  -1.16%  "
Set-TextValue $ws.Range("D15") "0.06963"
$ws.Range("E15").Value = "  -0.19%  "
Set-TextValue $ws.Range("D16") "88.17"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("E17").Value = "  -0.42%  "
Set-TextValue $ws.Range("D18") "0.000009966"
$ws.Range("E18").Value = "  -0.84%  "
Set-TextValue $ws.Range("D19") "16.96"
$ws.Range("E19").Value = "  -1.52%  "
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").Value = "28.941.46"
$ws.Range("E21").Value = "  +0.80%  "
Set-TextValue $ws.Range("D22") "5.245"
$ws.Range("E22").Value = "  -1.87%  "
$ws.Range("E23").Value = "  -0.60%  "
Set-TextValue $ws.Range("D24") "2.098"
$ws.Range("E24").Value = "  +2.06%  "
Set-TextValue $ws.Range("D25") "156.39"
$ws.Range("E25").Value = "  +0.90%  "
Set-TextValue $ws.Range("D26") "19.35"
$ws.Range("E26").Value = "  -1.72%  "
Set-TextValue $ws.Range("D27") "5.986"
$ws.Range("E27").Value = "  +2.60%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D28") "1.926"
$ws.Range("E28").Value = "  -2.01%  "
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D29") "117.55"
$ws.Range("E29").Value = "  -0.85%  "
Set-TextValue $ws.Range("D30") "0.09373"
$ws.Range("E30").Value = "  +0.23%  "
Set-TextValue $ws.Range("D31") "0.9037"
$ws.Range("E31").Value = "  -2.46%  "
Set-TextValue $ws.Range("D32") "5.269"
$ws.Range("E32").Value = "  -1.04%  "
Set-TextValue $ws.Range("D33") "1.319"
$ws.Range("E33").Value = "  -1.42%  "
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("E35").Value = "  +1.54%  "
Set-TextValue $ws.Range("D36") "0.05744"
$ws.Range("E36").Value = "  -0.71%  "
Set-TextValue $ws.Range("D37") "0.02075"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("E38").Value = "  -0.52%  "
Set-TextValue $ws.Range("D39") "7.658"
$ws.Range("E39").Value = "  -5.44%  "
Set-TextValue $ws.Range("D40") "0.5653"
$ws.Range("E40").Value = "  -1.48%  "
Set-TextValue $ws.Range("D41") "0.1764"
$ws.Range("E41").Value = "  -2.03%  "
Set-TextValue $ws.Range("D42") "9.699"
$ws.Range("E42").Value = "  -1.36%  "
Set-TextValue $ws.Range("D43") "2.272"
$ws.Range("E43").Value = "  +5.29%  "
Set-TextValue $ws.Range("D44") "11.91"
$ws.Range("E44").Value = "  +1.11%  "
Set-TextValue $ws.Range("D45") "0.5344"
$ws.Range("E45").Value = "  -1.08%  "
Set-TextValue $ws.Range("D46") "0.07046"
$ws.Range("E46").Value = "  -1.48%  "
Set-TextValue $ws.Range("D47") "1.846"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D48") "2.541"
$ws.Range("E48").Value = "  +2.19%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D49") "112.87"
$ws.Range("E49").Value = "  +0.30%  "
Set-TextValue $ws.Range("D50") "1.069"
$ws.Range("E50").Value = "  -4.56%  "
Set-TextValue $ws.Range("D51") "70.87"
$ws.Range("E51").Value = "  -0.33%  "
